# "fix plots for cases where aggregation happens in app"
#
# Adds a new issue row (question_id 1145) describing an R date-parsing
# error, and updates the worksheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row 10: question_id = 1145, issue_details = new note about a
# date-parsing failure ("in April (Month 04) 2018" form not recognized).
$ws.Range("C10").Value = 1145
$ws.Range("F10").Value = "R error; failure to parse date, date had form ""in April (Month 04) 2018"""

# Move/save the active selection to C11 (was H13).
$ws.Range("C11").Select() | Out-Null

# Best-effort: reposition the document window (xWindow/yWindow in the
# saved bookViews). Some hosts don't persist window chrome geometry back
# to OOXML, but setting it is harmless if unsupported.
$win = $excel.ActiveWindow
if ($win) {
    $win.Left = 16600
    $win.Top = 7800
}
